$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.229.48'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '1.862.60'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7046'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08167'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3034'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08168'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '1.843.55'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.171'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7088'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").Value = '29.239.64'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007894'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.42%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.785'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = '2.107.92'
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.413'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.961'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1442'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.10'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.967'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.428'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.486'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.386'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.057'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05205'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.170'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7075'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9981'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.668'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01848'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.733'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.17%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.143.04'
$ws.Range("E41").Value = '  +5.75%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9215'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4282'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.872'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9992'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.772'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("D49").Value = '2.004.71'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.197'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.955'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.17%  '
